$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.94
$ws.Range("E2").Value = 11.07

$ws.Range("B3").Value = 10.06
$ws.Range("E3").Value = 10.79

$ws.Range("E4").Value = 10.6
$ws.Range("F4").Value = 9.68
$ws.Range("H4").Value = 8.67

$ws.Range("B5").Value = 8.93
$ws.Range("C5").Value = 9.16
$ws.Range("D5").Value = 9.4
$ws.Range("F5").Value = 10.17
$ws.Range("G5").Value = 9.75

$ws.Range("D6").Value = 10.32
$ws.Range("E6").Value = 9.83
$ws.Range("G6").Value = 10.49
$ws.Range("H6").Value = 10.57

$ws.Range("E7").Value = 10.25
$ws.Range("F7").Value = 9.51
$ws.Range("H7").Value = 10.06

$ws.Range("D8").Value = 11.33
$ws.Range("F8").Value = 9.43
$ws.Range("G8").Value = 9.94
$ws.Range("J8").Value = 11.54

$ws.Range("H10").Value = 8.46
